# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns across the per-job Leve-profit sheets with the latest
# Universalis snapshot values.

$wb = $excel.ActiveWorkbook

# ---- ALC -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4850.2
$ws.Range("J51").Value = 4850.2
$ws.Range("L51").Value = 4850.2
$ws.Range("N51").Value = -5818.2

$ws.Range("H52").Value = 1981
$ws.Range("J52").Value = 2333
$ws.Range("L52").Value = 6999
$ws.Range("N52").Value = -7319

$ws.Range("H138").Value = 2250.946
$ws.Range("J138").Value = 2273.2122
$ws.Range("L138").Value = 6819.6366
$ws.Range("N138").Value = -17099.6366

# ---- ARM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2700.1667
$ws.Range("I2").Value = 1784
$ws.Range("J2").Value = 4532.5
$ws.Range("K2").Value = 1784
$ws.Range("L2").Value = 4532.5
$ws.Range("M2").Value = -1671
$ws.Range("N2").Value = -4758.5

$ws.Range("H45").Value = 1911.6957
$ws.Range("I45").Value = 1404.6666
$ws.Range("K45").Value = 1404.6666
$ws.Range("M45").Value = -1027.6666

$ws.Range("H97").Value = 479.2857
$ws.Range("J97").Value = 528
$ws.Range("L97").Value = 528
$ws.Range("N97").Value = -1520

$ws.Range("H116").Value = 2700.1667
$ws.Range("I116").Value = 1784
$ws.Range("J116").Value = 4532.5
$ws.Range("K116").Value = 1784
$ws.Range("L116").Value = 4532.5
$ws.Range("M116").Value = 510
$ws.Range("N116").Value = -9120.5

# ---- BSM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2700.1667
$ws.Range("I3").Value = 1784
$ws.Range("J3").Value = 4532.5
$ws.Range("K3").Value = 1784
$ws.Range("L3").Value = 4532.5
$ws.Range("M3").Value = -1670
$ws.Range("N3").Value = -4760.5

$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("M64").Value = -4775

$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("M67").Value = -4220

$ws.Range("H86").Value = 2172.1924
$ws.Range("I86").Value = 2122.1
$ws.Range("J86").Value = 2339.1667
$ws.Range("K86").Value = 2122.1
$ws.Range("L86").Value = 2339.1667
$ws.Range("M86").Value = -999.0999999999999
$ws.Range("N86").Value = -4585.1667

$ws.Range("H89").Value = 2172.1924
$ws.Range("I89").Value = 2122.1
$ws.Range("J89").Value = 2339.1667
$ws.Range("K89").Value = 10610.5
$ws.Range("L89").Value = 11695.8335
$ws.Range("M89").Value = -4994.5
$ws.Range("N89").Value = -22927.8335

$ws.Range("H94").Value = 1509.9231
$ws.Range("I94").Value = 536.6667
$ws.Range("K94").Value = 536.6667
$ws.Range("M94").Value = -85.66669999999999

# ---- CRP -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2246.1843
$ws.Range("I99").Value = 2219.9119
$ws.Range("J99").Value = 2469.5
$ws.Range("K99").Value = 2219.9119
$ws.Range("L99").Value = 2469.5
$ws.Range("M99").Value = -721.9119000000001
$ws.Range("N99").Value = -5465.5

$ws.Range("H126").Value = 2246.1843
$ws.Range("I126").Value = 2219.9119
$ws.Range("J126").Value = 2469.5
$ws.Range("K126").Value = 6659.7357
$ws.Range("L126").Value = 7408.5
$ws.Range("M126").Value = -4189.7357
$ws.Range("N126").Value = -12348.5

# ---- CUL -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 999
$ws.Range("J22").Value = 999
$ws.Range("L22").Value = 2997
$ws.Range("N22").Value = -3335

$ws.Range("H23").Value = 9653332
$ws.Range("I23").Value = 469.66666
$ws.Range("J23").Value = 14479764
$ws.Range("K23").Value = 1408.99998
$ws.Range("L23").Value = 43439292
$ws.Range("M23").Value = -1173.99998
$ws.Range("N23").Value = -43439762

$ws.Range("H27").Value = 999
$ws.Range("J27").Value = 999
$ws.Range("L27").Value = 2997
$ws.Range("N27").Value = -3201

$ws.Range("H35").Value = 501
$ws.Range("I35").Value = 99
$ws.Range("K35").Value = 297
$ws.Range("M35").Value = -9

$ws.Range("H41").Value = 1009.8333
$ws.Range("I41").Value = 20
$ws.Range("K41").Value = 60
$ws.Range("M41").Value = 278

$ws.Range("H59").Value = 1001
$ws.Range("I59").Value = 1001
$ws.Range("K59").Value = 3003
$ws.Range("M59").Value = -2463

$ws.Range("H69").Value = 3955.7144
$ws.Range("I69").Value = 4245
$ws.Range("J69").Value = 2220
$ws.Range("K69").Value = 12735
$ws.Range("L69").Value = 6660
$ws.Range("M69").Value = -11924
$ws.Range("N69").Value = -8282

$ws.Range("H72").Value = 3955.7144
$ws.Range("I72").Value = 4245
$ws.Range("J72").Value = 2220
$ws.Range("K72").Value = 38205
$ws.Range("L72").Value = 19980
$ws.Range("M72").Value = -34149
$ws.Range("N72").Value = -28092

# ---- GSM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3210.75
$ws.Range("I102").Value = 3091.625
$ws.Range("J102").Value = 3449
$ws.Range("K102").Value = 3091.625
$ws.Range("L102").Value = 3449
$ws.Range("M102").Value = -1469.625
$ws.Range("N102").Value = -6693

# ---- LTW -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5173.0435
$ws.Range("I7").Value = 3771.2307
$ws.Range("K7").Value = 3771.2307
$ws.Range("M7").Value = -3659.2307

$ws.Range("H40").Value = 6797.1577
$ws.Range("I40").Value = 6420.5293
$ws.Range("K40").Value = 6420.5293
$ws.Range("M40").Value = -6284.5293

$ws.Range("H55").Value = 8443.76
$ws.Range("J55").Value = 23080.334
$ws.Range("L55").Value = 23080.334
$ws.Range("N55").Value = -23426.334

$ws.Range("H93").Value = 5208.273
$ws.Range("I93").Value = 6881.5
$ws.Range("K93").Value = 6881.5
$ws.Range("M93").Value = -5633.5

$ws.Range("H126").Value = 5173.0435
$ws.Range("I126").Value = 3771.2307
$ws.Range("K126").Value = 11313.6921
$ws.Range("M126").Value = -8843.6921

# ---- WVR -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 111114960
$ws.Range("I126").Value = 3801.4
$ws.Range("K126").Value = 11404.2
$ws.Range("M126").Value = -8934.200000000001
